$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 254, shifting old rows 254-261
# down to 257-264, then populate the 3 new rows with the new "Clementina"
# price records.
$ws.Rows("254:256").Insert()

$rows = @(254, 255, 256)
$data = @(
    @{ L = "Especial"; M = 40; N = 14000; O = 14000; P = 14000; S = 1400 },
    @{ L = "Primera";  M = 50; N = 12000; O = 12000; P = 12000; S = 1200 },
    @{ L = "Segunda";  M = 30; N = 10000; O = 10000; P = 10000; S = 1000 }
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $d = $data[$i]

    $ws.Cells.Item($r, 1).Value = 7
    $ws.Cells.Item($r, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($r, 3).Value = "Ñuble"
    $ws.Cells.Item($r, 4).Value = 45075
    $ws.Cells.Item($r, 5).Value = 16
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = "Cítricos"
    $ws.Cells.Item($r, 9).Value = 100102004
    $ws.Cells.Item($r, 10).Value = "Mandarina"
    $ws.Cells.Item($r, 11).Value = "Clementina"
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
    $ws.Cells.Item($r, 14).Value = $d.N
    $ws.Cells.Item($r, 15).Value = $d.O
    $ws.Cells.Item($r, 16).Value = $d.P
    $ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"
    $ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($r, 19).Value = $d.S
    $ws.Cells.Item($r, 20).Value = 10
}
